$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @(0.4145870963266276,0.1615748415685587,0.1784203654063781,1.907250001524517,1.258863221492163,1.139205889773478,0.2336961426485686,0.1609786650763496,4.892343549182613)
  3 = @(0.4107024072650205,0.1615709920222912,0.1778979785889234,1.877333202784769,1.226370776659252,1.12911278080685,0.2324819634782926,0.1611842094933174,4.801883612414315)
  4 = @(0.4085160041098703,0.1616226682646627,0.1776499131432097,1.859969350165528,1.207147119791188,1.12346002733662,0.2318405107042523,0.1613630788569225,4.749003842339505)
  5 = @(0.4076750433128353,0.1616573420265084,0.1775671208319984,1.853146138553939,1.19949587067164,1.121293439370419,0.2316052930080659,0.1614492342832712,4.728124154461455)
  6 = @(0.407538423988143,0.1616639226903445,0.1775544788789674,1.852028410230815,1.19823640563169,1.120941952292839,0.2315678168708715,0.1614643421286104,4.724697508323914)
  7 = @(0.4085044600716969,0.1616230807243824,0.1776487224690619,1.859876306866084,1.20704319353959,1.123430253392598,0.2318372324579769,0.1613641870419684,4.748719542336403)
  8 = @(0.4132064037992222,0.1615623138245326,0.1782251725284851,1.896726010695119,1.24750882414645,1.135612805412848,0.2332559032986694,0.1610386191530893,4.860599982157794)
  9 = @(0.4240043641165983,0.1618709033210948,0.1799316716273083,1.97697307805548,1.332644793882451,1.163823441513983,0.2368634153450557,0.1608171973375399,5.10116252340174)
  10 = @(0.4329007671991292,0.1623572111732869,0.1815361776646114,2.040820501381802,1.3987501284129,1.187188775875427,0.2400175384908252,0.1609077591747052,5.29088259680185)
  11 = @(0.437157536141342,0.1626345781304366,0.1823422034695525,2.070933559428283,1.429602432515935,1.198392676901079,0.2415619452893978,0.1610037532331106,5.380027571598816)
  12 = @(0.4387996321026719,0.1627476604520695,0.1826583594943827,2.08249050447148,1.44139806064598,1.202718009387411,0.2421625299240802,0.1610479651501322,5.414193847512422)
  13 = @(0.4384446366102281,0.1627229485252997,0.1825897836986812,2.079994667914917,1.438852651878591,1.201782795529311,0.2420324828275469,0.1610380939956855,5.406817328695979)
  14 = @(0.4372920281788311,0.1626437202936017,0.1823679948194865,2.071881273450202,1.430570609092314,1.198746868062045,0.2416110401033933,0.1610072331328709,5.382830251232178)
  15 = @(0.4365899483082671,0.1625962382890407,0.1822335658805763,2.066931616112072,1.425512281660389,1.196898039116945,0.2413549453743116,0.1609893531540649,5.36819075578677)
  16 = @(0.4326267970164395,0.1623402122907649,0.1814850325341304,2.038874054898656,1.396749589797309,1.186468141939258,0.2399188123671863,0.1609025867895042,5.285113966609515)
  17 = @(0.4302492454531262,0.162197511197526,0.1810453185918952,2.021935419282357,1.379304747435526,1.180216973275918,0.2390658528427707,0.1608633798888732,5.234876931550048)
  18 = @(0.4289014840590539,0.1621207191679446,0.180799573198211,2.012293332637228,1.369344397988584,1.17667557853423,0.2385855676478457,0.1608459879751223,5.206249231189076)
  19 = @(0.4284485465299497,0.1620956272873002,0.1807175993517802,2.009045956693328,1.365984603809324,1.175485817301706,0.2384247230671477,0.1608409859847058,5.19660227799335)
  20 = @(0.4305002966017355,0.162212155078997,0.1810913853069458,2.023728156611398,1.381154175649641,1.180876820556904,0.2391555843693425,0.1608670197441953,5.240197079765721)
  21 = @(0.4376297590447109,0.1626667732733722,0.1824328430330162,2.074260199388647,1.433000189587574,1.199636350297197,0.2417344006087063,0.1610160845177546,5.38986473235741)
  22 = @(0.4424650069443601,0.1630108001391619,0.1833732736285434,2.108182241695914,1.467540598247581,1.212378580641257,0.243511614216203,0.1611593186795943,5.49006548527683)
  23 = @(0.4398682683589925,0.1628229022961492,0.1828655236690686,2.089995342099598,1.449045613039544,1.20553373112412,0.2425546833322372,0.1610786862495246,5.436368103370228)
  24 = @(0.4303867368065539,0.1622055182295625,0.1810705365783676,2.022917360543076,1.380317834488039,1.180578340109975,0.239114985309314,0.1608653581272819,5.237791048862846)
  25 = @(0.4209142018315504,0.1617417508112524,0.1794083929520234,1.954407016463563,1.308991014268571,1.155728654418482,0.2357990689262337,0.1608325559549506,5.033810303472819)
}

$cols = @("C", "D", "E", "F", "G", "H", "J", "L", "O")
foreach ($r in $data.Keys) {
  $vals = $data[$r]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + $r).Value = $vals[$i]
  }
}

Write-Output "done"